$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.468507333333333
$ws.Range("H2").Value = 4.405521999999999
$ws.Range("I2").Value = 0.005118279455112885
$ws.Range("J2").Value = 0.005118279455112885
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02725333333333333
$ws.Range("N2").Value = 0.08176
$ws.Range("O2").Value = 0.0007089206372884383
$ws.Range("P2").Value = 0.0007089206372884382
$ws.Range("Q2").Value = 0.04002171985777778
$ws.Range("R2").Value = 0.36019547872
$ws.Range("S2").Value = 0.000003628453933138947
$ws.Range("T2").Value = 0.000003628453933138947

$ws.Range("G3").Value = 1.468507333333333
$ws.Range("H3").Value = 4.405521999999999
$ws.Range("I3").Value = 0.005118279455112885
$ws.Range("J3").Value = 0.005118279455112885
$ws.Range("M3").Value = 38.416166
$ws.Range("N3").Value = 115.248498
$ws.Range("O3").Value = 0.9992910793627116
$ws.Range("P3").Value = 0.9992910793627116
$ws.Range("Q3").Value = 56.41442148955066
$ws.Range("R3").Value = 507.7297934059559
$ws.Range("S3").Value = 0.005114651001179746
$ws.Range("T3").Value = 0.005114651001179746

$ws.Range("I4").Value = 0.9046276674881553
$ws.Range("J4").Value = 0.9046276674881553
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02725333333333333
$ws.Range("N4").Value = 0.08176
$ws.Range("O4").Value = 0.0007089206372884383
$ws.Range("P4").Value = 0.0007089206372884382
$ws.Range("Q4").Value = 7.073618273742223
$ws.Range("R4").Value = 63.66256446368001
$ws.Range("S4").Value = 0.0006413092225444565
$ws.Range("T4").Value = 0.0006413092225444565

$ws.Range("I5").Value = 0.9046276674881553
$ws.Range("J5").Value = 0.9046276674881553
$ws.Range("M5").Value = 38.416166
$ws.Range("N5").Value = 115.248498
$ws.Range("O5").Value = 0.9992910793627116
$ws.Range("P5").Value = 0.9992910793627116
$ws.Range("Q5").Value = 9970.93788495773
$ws.Range("R5").Value = 89738.44096461957
$ws.Range("S5").Value = 0.9039863582656108
$ws.Range("T5").Value = 0.9039863582656108

$ws.Range("G6").Value = 0.5890733333333333
$ws.Range("H6").Value = 1.76722
$ws.Range("I6").Value = 0.002053133730501083
$ws.Range("J6").Value = 0.002053133730501083
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02725333333333333
$ws.Range("N6").Value = 0.08176
$ws.Range("O6").Value = 0.0007089206372884383
$ws.Range("P6").Value = 0.0007089206372884382
$ws.Range("Q6").Value = 0.01605421191111111
$ws.Range("R6").Value = 0.1444879072
$ws.Range("S6").Value = 0.000001455508872665217
$ws.Range("T6").Value = 0.000001455508872665216

$ws.Range("G7").Value = 0.5890733333333333
$ws.Range("H7").Value = 1.76722
$ws.Range("I7").Value = 0.002053133730501083
$ws.Range("J7").Value = 0.002053133730501083
$ws.Range("M7").Value = 38.416166
$ws.Range("N7").Value = 115.248498
$ws.Range("O7").Value = 0.9992910793627116
$ws.Range("P7").Value = 0.9992910793627116
$ws.Range("Q7").Value = 22.62993895950666
$ws.Range("R7").Value = 203.66945063556
$ws.Range("S7").Value = 0.002051678221628418
$ws.Range("T7").Value = 0.002051678221628418

$ws.Range("G8").Value = 25.306101
$ws.Range("H8").Value = 75.918303
$ws.Range("I8").Value = 0.0882009193262308
$ws.Range("J8").Value = 0.0882009193262308
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02725333333333333
$ws.Range("N8").Value = 0.08176
$ws.Range("O8").Value = 0.0007089206372884383
$ws.Range("P8").Value = 0.0007089206372884382
$ws.Range("Q8").Value = 0.68967560592
$ws.Range("R8").Value = 6.20708045328
$ws.Range("S8").Value = 0.00006252745193817767
$ws.Range("T8").Value = 0.00006252745193817766

$ws.Range("G9").Value = 25.306101
$ws.Range("H9").Value = 75.918303
$ws.Range("I9").Value = 0.0882009193262308
$ws.Range("J9").Value = 0.0882009193262308
$ws.Range("M9").Value = 38.416166
$ws.Range("N9").Value = 115.248498
$ws.Range("O9").Value = 0.9992910793627116
$ws.Range("P9").Value = 0.9992910793627116
$ws.Range("Q9").Value = 972.1633768287659
$ws.Range("R9").Value = 8749.470391458894
$ws.Range("S9").Value = 0.08813839187429262
$ws.Range("T9").Value = 0.08813839187429262
